$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.550.56'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.922.47'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.68'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4876'
$ws.Range('E7').Value = '  +3.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2905'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06718'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '110.79'
$ws.Range('E10').Value = '  +5.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.22'
$ws.Range('E11').Value = '  +4.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.921.56'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07584'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.339'
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6712'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '294.99'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.592.69'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.06'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.561'
$ws.Range('E20').Value = '  +2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007572'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.177.71'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.491'
$ws.Range('E24').Value = '  +2.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.478'
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.52'
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.23'
$ws.Range('E27').Value = '  -3.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.114'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1074'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  +5.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.163'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.053'
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05049'
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7396'
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.140'
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9993'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.713'
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02030'
$ws.Range('E38').Value = '  -2.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.686'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.70'
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.026'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4437'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8651'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.09'
$ws.Range('E44').Value = '  +5.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.824'
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.250'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.39'
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.265'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1231'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  +4.53%  '
